$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rule R40 (row 11) now reports the literal text "1" instead of "R40".
# Leading apostrophe forces the numeric-looking value to be stored as text
# (matches the sharedStrings <t>1</t> entry / t="s" cell in the target file).
$ws.Range("B11").Value = "'1"
